# feat: add 2022-Q4 data
#
# 1) Insert a new "2022-Q4" row at the top of the "总计" (totals) summary
#    sheet, pushing the existing Q3/Q2/Q1 rows down by one.
# 2) Insert a brand-new "2022-Q4" worksheet (with the per-fund holdings
#    detail) right after "总计", before the existing "2022-Q3" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: shift rows 2-4 down to 3-5, write new row 2 = 2022-Q4
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Clone the formatting of the last existing data row onto the new row 5
# first (so the freshly-used A5 cell picks up the same index-column
# style as A2:A4), THEN overwrite the values top-to-bottom.
$totals.Range("A4").Copy($totals.Range("A5"))

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 3
$totals.Range("D2").Value = 0.25

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("C3").Value = 18
$totals.Range("D3").Value = 5.89

$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2022-Q2"
$totals.Range("C4").Value = 15
$totals.Range("D4").Value = 5.17

$totals.Range("A5").Value = 3
$totals.Range("B5").Value = "2022-Q1"
$totals.Range("C5").Value = 16
$totals.Range("D5").Value = 7.89

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet, inserted right after "总计"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item("总计"))
$q4.Name = "2022-Q4"

# Match the outline defaults used by the sibling quarter sheets.
$q4.Outline.SummaryRow = 1
$q4.Outline.SummaryColumn = 1
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# Header row, copying the "2022-Q3" header style (bold, bordered) across.
# (Re-fetch the sheet reference fresh, rather than reusing one captured
# before the worksheet insertion above shifted sheet indices.)
$wb.Worksheets.Item("2022-Q3").Range("B1:H1").Copy($q4.Range("B1:H1"))
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Index column (A) styling, matching the sibling sheets' "s=2" style.
$wb.Worksheets.Item("2022-Q3").Range("A2").Copy($q4.Range("A2"))
$wb.Worksheets.Item("2022-Q3").Range("A2").Copy($q4.Range("A3"))
$wb.Worksheets.Item("2022-Q3").Range("A2").Copy($q4.Range("A4"))

# Columns B-G hold text in the source data (fund codes, names, and
# numeric-looking figures alike are all stored as text), so force the
# Text number format before writing the values.
$q4.Range("B2:G4").NumberFormat = "@"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "513300"
$q4.Range("C2").Value = "华夏纳斯达克100ETF（QDII）"
$q4.Range("D2").Value = "12.43"
$q4.Range("E2").Value = "97.54"
$q4.Range("F2").Value = "1.81"
$q4.Range("G2").Value = "0.2250"
$q4.Range("H2").Value = 9

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "013328"
$q4.Range("C3").Value = "嘉实全球价值股票（QDII）人民币"
$q4.Range("D3").Value = "1.62"
$q4.Range("E3").Value = "90.66"
$q4.Range("F3").Value = "0.79"
$q4.Range("G3").Value = "0.0128"
$q4.Range("H3").Value = 9

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "013329"
$q4.Range("C4").Value = "嘉实全球价值股票（QDII）美元现汇"
$q4.Range("D4").Value = "1.62"
$q4.Range("E4").Value = "90.66"
$q4.Range("F4").Value = "0.79"
$q4.Range("G4").Value = "0.0128"
$q4.Range("H4").Value = 9
